$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 1.08
$ws.Range("K3").Value = 8
$ws.Range("L3").Value = 1.44
$ws.Range("M3").Value = 2.63
$ws.Range("N3").Value = 2.35
$ws.Range("O3").Value = 1.57
$ws.Range("R3").Value = 2.2
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 5.5
$ws.Range("U3").Value = 7
$ws.Range("Z3").Value = 7.5
$ws.Range("AB3").Value = 21
$ws.Range("AC3").Value = 81
$ws.Range("AE3").Value = 10

# Row 6
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 3.25
$ws.Range("N6").Value = 2.05
$ws.Range("O6").Value = 1.75
$ws.Range("R6").Value = 1.91
$ws.Range("S6").Value = 1.8
$ws.Range("T6").Value = 6.5
$ws.Range("V6").Value = 8.5
$ws.Range("Y6").Value = 29
$ws.Range("Z6").Value = 9
$ws.Range("AA6").Value = 6.5
$ws.Range("AB6").Value = 17
$ws.Range("AC6").Value = 51
$ws.Range("AE6").Value = 12
$ws.Range("AJ6").Value = 41

# Row 9
$ws.Range("C9").Value = "19:00"
$ws.Range("K9").Value = 13
$ws.Range("L9").Value = 1.22
$ws.Range("M9").Value = 4

# Row 10
$ws.Range("I10").Value = 8.5

# Row 11
$ws.Range("G11").Value = 1.9
$ws.Range("H11").Value = 3.3
$ws.Range("J11").Value = 1.08
$ws.Range("K11").Value = 8
$ws.Range("L11").Value = 1.44
$ws.Range("M11").Value = 2.63
$ws.Range("N11").Value = 2.4
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 1.53
$ws.Range("Q11").Value = 2.38
$ws.Range("X11").Value = 19
$ws.Range("AE11").Value = 9.5
$ws.Range("AJ11").Value = 51

# Row 12
$ws.Range("I12").Value = 6.2
$ws.Range("L12").Value = 1.37
$ws.Range("M12").Value = 2.62
$ws.Range("N12").Value = 2.07
$ws.Range("O12").Value = 1.6
$ws.Range("Q12").Value = 2.37
$ws.Range("R12").Value = 2.22
$ws.Range("S12").Value = 1.52
$ws.Range("T12").Value = 5.2
$ws.Range("U12").Value = 5.9
$ws.Range("Z12").Value = 8
$ws.Range("AB12").Value = 25
$ws.Range("AC12").Value = 175
$ws.Range("AE12").Value = 12.5
$ws.Range("AG12").Value = 22
$ws.Range("AI12").Value = 90

# Row 15
$ws.Range("G15").Value = 4.45
$ws.Range("H15").Value = 3.8
$ws.Range("I15").Value = 1.65
$ws.Range("L15").Value = 1.19
$ws.Range("M15").Value = 3.7
$ws.Range("N15").Value = 1.55
$ws.Range("O15").Value = 2.15
$ws.Range("R15").Value = 1.55
$ws.Range("S15").Value = 2.15
$ws.Range("T15").Value = 17
$ws.Range("U15").Value = 30
$ws.Range("V15").Value = 14
$ws.Range("W15").Value = 75
$ws.Range("X15").Value = 37
$ws.Range("Y15").Value = 32
$ws.Range("Z15").Value = 14.5
$ws.Range("AA15").Value = 7.8
$ws.Range("AB15").Value = 13
$ws.Range("AC15").Value = 45
$ws.Range("AE15").Value = 8.75
$ws.Range("AF15").Value = 9.25
$ws.Range("AG15").Value = 8
$ws.Range("AH15").Value = 14
$ws.Range("AI15").Value = 12
$ws.Range("AJ15").Value = 20

# Row 16
$ws.Range("I16").Value = 3.6
$ws.Range("V16").Value = 8.5
$ws.Range("AE16").Value = 15

# Row 17
$ws.Range("G17").Value = 1.33
$ws.Range("H17").Value = 5.25
$ws.Range("I17").Value = 9.5
$ws.Range("L17").Value = 1.18
$ws.Range("M17").Value = 4.5
$ws.Range("N17").Value = 1.6
$ws.Range("O17").Value = 2.3
$ws.Range("R17").Value = 2
$ws.Range("S17").Value = 1.75
$ws.Range("U17").Value = 6.5
$ws.Range("V17").Value = 9
$ws.Range("W17").Value = 8.5
$ws.Range("Y17").Value = 29
$ws.Range("Z17").Value = 15
$ws.Range("AA17").Value = 10
$ws.Range("AB17").Value = 21
$ws.Range("AC17").Value = 67
$ws.Range("AD17").Value = 351
$ws.Range("AG17").Value = 23
$ws.Range("AH17").Value = 101

# Row 20
$ws.Range("K20").Value = 29
$ws.Range("U20").Value = 9
$ws.Range("AF20").Value = 67

# Row 21
$ws.Range("N21").Value = 1.62
$ws.Range("O21").Value = 2.25

# Row 23
$ws.Range("G23").Value = 2.22
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 3.25
$ws.Range("L23").Value = 1.44
$ws.Range("M23").Value = 2.42
$ws.Range("N23").Value = 2.25
$ws.Range("O23").Value = 1.5
$ws.Range("R23").Value = 1.93
$ws.Range("S23").Value = 1.7
$ws.Range("T23").Value = 6.3
$ws.Range("U23").Value = 9.75
$ws.Range("V23").Value = 9.25
$ws.Range("W23").Value = 22
$ws.Range("X23").Value = 21
$ws.Range("Z23").Value = 7
$ws.Range("AA23").Value = 5.9
$ws.Range("AB23").Value = 17
$ws.Range("AC23").Value = 100
$ws.Range("AE23").Value = 7.7
$ws.Range("AF23").Value = 15.5
$ws.Range("AG23").Value = 12
$ws.Range("AH23").Value = 45
$ws.Range("AI23").Value = 35
$ws.Range("AJ23").Value = 50
